# Generate Report for handoff
# - Overview sheet: status text "Handoff transform failed" -> "Ready for handoff"
# - zh-cn / de-de sheets: same status text change, plus newly generated
#   handoff file hyperlink (column C), handoff datetime (column D), and
#   handoff reason changed from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/212e4a0128a0aa603c3efb97be5207021e16a30e/e2e/e6f9569e-0936-4cc8-86d0-5368f631f522.809576b6647759f740283803b33c99d09161c5d8.zh-cn.xlf",
    "",
    "",
    "e6f9569e-0936-4cc8-86d0-5368f631f522.809576b6647759f740283803b33c99d09161c5d8.zh-cn.xlf"
) | Out-Null
$wsZhCn.Range("D2").Value = "2016-01-14 03:37:01"
$wsZhCn.Range("H2").Value = "Include"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = "Ready for handoff"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/212e4a0128a0aa603c3efb97be5207021e16a30e/e2e/e6f9569e-0936-4cc8-86d0-5368f631f522.809576b6647759f740283803b33c99d09161c5d8.de-de.xlf",
    "",
    "",
    "e6f9569e-0936-4cc8-86d0-5368f631f522.809576b6647759f740283803b33c99d09161c5d8.de-de.xlf"
) | Out-Null
$wsDeDe.Range("D2").Value = "2016-01-14 03:37:12"
$wsDeDe.Range("H2").Value = "Include"
